$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows whose text changed, and append the new rows for the
# "AutoPlay/Pause", "Doppelbesetzung" and "Spiel mit zwei Menschen" test cases.

$ws.Range("A8").Value = "Der Nutzer klickt auf das `"NextMove`"-Symbol."
$ws.Range("B8").Value = "Die KI macht einen zufälligen Zug und belegt somit ein zufälliges unbesetztes Feld. Dieses wird mit dem Zeichen `"O`" versehen. Im Graph wird der erneuerte Verlauf plus alle möglichen Folgezustände angezeigt."

$ws.Range("A14").Value = "Der Nutzer klickt auf das `"Play`"-Symbol."
$ws.Range("B14").Value = "Das `"Play`"-Symbol wird durch ein `"Pause`"-Symbol ersetzt. Es wird kein Zug ausgeführt, da der Mensch am Zug ist."

$ws.Range("A15").Value = "Der Nutzer klickt auf das `"NextMove`"-Symbol."
$ws.Range("B15").Value = "Das `"Pause`"-Symbol wird durch ein `"Play`"-Symbol ersetzt. Es wird kein Zug ausgeführt, da der Mensch am Zug ist."

$ws.Range("A16").Value = "Der Nutzer klickt auf ein Feld des Spielfeldes."
$ws.Range("B16").Value = "Das angeklickte Feld wird mit dem Zeichen `"X`" versehen. Im Graph wird jetzt der Verlauf plus alle weiteren möglichen Folgezustände angezeigt."

$ws.Range("A17").Value = "Der Nutzer klickt auf das `"NextMove`"-Symbol."
$ws.Range("B17").Value = "Die KI macht einen zufälligen Zug und belegt somit ein zufälliges unbesetztes Feld. Dieses wird mit dem Zeichen `"O`" versehen. Im Graph wird der erneuerte Verlauf plus alle möglichen Folgezustände angezeigt."

$ws.Range("A18").Value = "Der Nutzer klickt auf ein unbesetztes Feld des Spielfeldes."
$ws.Range("B18").Value = "Das angeklickte Feld wird mit dem Zeichen `"X`" versehen. Im Graph wird jetzt der Verlauf plus alle weiteren möglichen Folgezustände angezeigt."

$ws.Range("A19").Value = "Der Nutzer klickt auf das `"Play`"-Symbol."
$ws.Range("B19").Value = "Das `"Play`"-Symbol wird durch ein `"Pause`"-Symbol ersetzt. Die KI macht einen zufälligen Zug und belegt somit ein zufälliges unbesetztes Feld. Dieses wird mit dem Zeichen `"O`" versehen. Im Graph wird der erneuerte Verlauf plus alle möglichen Folgezustände angezeigt."

$ws.Range("A20").Value = "Der Nutzer klickt auf ein unbesetztes Feld des Spielfeldes."
$ws.Range("B20").Value = "Das angeklickte Feld wird mit dem Zeichen `"X`" versehen. Im Graph wird jetzt der Verlauf plus alle weiteren möglichen Folgezustände angezeigt. Kurze Zeit später macht die KI einen zufälligen Zug und belegt somit ein zufälliges unbesetztes Feld. Dieses wird mit dem Zeichen `"O`" versehen. Im Graph wird der erneuerte Verlauf plus alle möglichen Folgezustände angezeigt."

$ws.Range("A21").Value = "Der Nutzer klickt auf das `"Pause`"-Symbol."
$ws.Range("B21").Value = "Das `"Pause`"-Symbol wird durch ein `"Play`"-Symbol ersetzt."

$ws.Range("A23").Value = "Test: Doppelbesetzung"
$ws.Range("A23").Font.Bold = $true

$ws.Range("A24").Value = "Der Nutzer klickt auf ein bereits besetztes Feld des Spielfeldes."
$ws.Range("B24").Value = "Es passiert nichts, da das Feld bereits besetzt ist."

$ws.Range("A25").Value = "Der Nutzer klickt auf ein unbesetztes Feld des Spielfeldes."
$ws.Range("B25").Value = "Das angeklickte Feld wird mit dem Zeichen `"X`" versehen. Im Graph wird jetzt der Verlauf plus alle weiteren möglichen Folgezustände angezeigt."

$ws.Range("A27").Value = "Test: Spiel mit zwei Menschen"
$ws.Range("A27").Font.Bold = $true

$ws.Range("A28").Value = "Der Nutzer klickt auf den Button `"Neustart`"."
$ws.Range("B28").Value = "Die Spielerauswahl wird angezeigt."

$ws.Range("A29").Value = "Der Nutzer klickt auf das DropDown-Menü für die Auswahl von Spieler 2."
$ws.Range("B29").Value = "Es wird im DropDown-Menü eine List aller möglichen Spieler angezeigt."

$ws.Range("A30").Value = "Der Nutzer wählt die Option `"Mensch`" im DropDown-Menü aus."
$ws.Range("B30").Value = "Als Spieler 2 wird ein Mensch festgelegt. Im DropDown-Menü wird angezeigt, dass die Option ausgewählt wurde."

$ws.Range("A31").Value = "Der Nutzer behält die Konfiguration bei und klickt auf `"Spiel starten`""
$ws.Range("B31").Value = "Es wird in die Spielansicht gewechselt. Das Spielfeld ist leer. Der Graph zeigt das leere Feld und alle möglichen Folgezustände an."

$ws.Range("A32").Value = "Der Nutzer klickt auf das `"Play`"-Symbol."
$ws.Range("B32").Value = "Das `"Play`"-Symbol wird durch ein `"Pause`"-Symbol ersetzt. Es wird kein Zug ausgeführt, da der Mensch am Zug ist."

$ws.Range("A33").Value = "Der Nutzer klickt auf das `"NextMove`"-Symbol."
$ws.Range("B33").Value = "Das `"Pause`"-Symbol wird durch ein `"Play`"-Symbol ersetzt. Es wird kein Zug ausgeführt, da der Mensch am Zug ist."

$ws.Range("A34").Value = "Der Nutzer klickt auf ein unbesetztes Feld des Spielfeldes."
$ws.Range("B34").Value = "Das angeklickte Feld wird mit dem Zeichen `"X`" versehen. Im Graph wird jetzt der Verlauf plus alle weiteren möglichen Folgezustände angezeigt."

$ws.Range("A35").Value = "Der Nutzer klickt auf ein unbesetztes Feld des Spielfeldes."
$ws.Range("B35").Value = "Das angeklickte Feld wird mit dem Zeichen `"O`" versehen. Im Graph wird jetzt der Verlauf plus alle weiteren möglichen Folgezustände angezeigt."

$ws.Range("A36").Value = "Der Nutzer klickt auf ein unbesetztes Feld des Spielfeldes."
$ws.Range("B36").Value = "Das angeklickte Feld wird mit dem Zeichen `"X`" versehen. Im Graph wird jetzt der Verlauf plus alle weiteren möglichen Folgezustände angezeigt."

# Move the active selection to A37, matching the post-edit cursor position.
$ws.Range("A37").Select()